$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N17").Value = -2489.625
$ws.Range("L17").Value = 2153.625
$ws.Range("H17").Value = 693.1177
$ws.Range("J17").Value = 717.875
$ws.Range("I100").Value = 850
$ws.Range("L100").Value = 1041
$ws.Range("K100").Value = 850
$ws.Range("M100").Value = -309
$ws.Range("N100").Value = -2123
$ws.Range("H100").Value = 1002.8
$ws.Range("J100").Value = 1041
$ws.Range("H101").Value = 552.6667
$ws.Range("J101").Value = 1492.5
$ws.Range("I101").Value = 284.14285
$ws.Range("L101").Value = 4477.5
$ws.Range("K101").Value = 852.4285500000001
$ws.Range("M101").Value = 769.5714499999999
$ws.Range("N101").Value = -7721.5
$ws.Range("N103").Value = -2573
$ws.Range("L103").Value = 1401
$ws.Range("K103").Value = 1398
$ws.Range("H103").Value = 466.89285
$ws.Range("M103").Value = -812
$ws.Range("J103").Value = 467
$ws.Range("I103").Value = 466
$ws.Range("I111").Value = 1675.4166
$ws.Range("L111").Value = 11157.75
$ws.Range("K111").Value = 5026.2498
$ws.Range("M111").Value = -1959.2498
$ws.Range("N111").Value = -17291.75
$ws.Range("H111").Value = 2492.95
$ws.Range("J111").Value = 3719.25
$ws.Range("M138").Value = -104.5712000000003
$ws.Range("N138").Value = -33174.905
$ws.Range("H138").Value = 5939.137
$ws.Range("J138").Value = 7631.635
$ws.Range("I138").Value = 1748.1904
$ws.Range("L138").Value = 22894.905
$ws.Range("K138").Value = 5244.5712

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L45").Value = 2057.8
$ws.Range("N45").Value = -2811.8
$ws.Range("H45").Value = 72792.93
$ws.Range("J45").Value = 2057.8
$ws.Range("J74").Value = 3502
$ws.Range("I74").Value = 605.6
$ws.Range("L74").Value = 3502
$ws.Range("K74").Value = 605.6
$ws.Range("M74").Value = 268.4
$ws.Range("N74").Value = -5250
$ws.Range("H74").Value = 1239.1875
$ws.Range("H77").Value = 1239.1875
$ws.Range("K77").Value = 3028
$ws.Range("M77").Value = 1340
$ws.Range("J77").Value = 3502
$ws.Range("I77").Value = 605.6
$ws.Range("N77").Value = -26246
$ws.Range("L77").Value = 17510
$ws.Range("K97").Value = 77853.695
$ws.Range("M97").Value = -77357.695
$ws.Range("N97").Value = -4418.6
$ws.Range("H97").Value = 57179.5
$ws.Range("J97").Value = 3426.6
$ws.Range("I97").Value = 77853.695
$ws.Range("L97").Value = 3426.6
$ws.Range("I102").Value = 144988.42
$ws.Range("L102").Value = 3095.8
$ws.Range("K102").Value = 144988.42
$ws.Range("M102").Value = -143366.42
$ws.Range("N102").Value = -6339.8
$ws.Range("H102").Value = 61522.176
$ws.Range("J102").Value = 3095.8
$ws.Range("I110").Value = 100200910
$ws.Range("L110").Value = 1535.2
$ws.Range("K110").Value = 100200910
$ws.Range("M110").Value = -100198865
$ws.Range("N110").Value = -5625.2
$ws.Range("H110").Value = 66801120
$ws.Range("J110").Value = 1535.2
$ws.Range("J122").Value = 1380
$ws.Range("I122").Value = 1333.6875
$ws.Range("N122").Value = -9040
$ws.Range("L122").Value = 4140
$ws.Range("M122").Value = -1551.0625
$ws.Range("K122").Value = 4001.0625
$ws.Range("H122").Value = 1347.7826
$ws.Range("H132").Value = 2932.3901
$ws.Range("I132").Value = 2922.3784
$ws.Range("M132").Value = -6237.135200000001
$ws.Range("K132").Value = 8767.1352

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K64").Value = 99.5
$ws.Range("M64").Value = 125.5
$ws.Range("N64").Value = -1296.75
$ws.Range("L64").Value = 846.75
$ws.Range("H64").Value = 597.6667
$ws.Range("J64").Value = 846.75
$ws.Range("I64").Value = 99.5
$ws.Range("I67").Value = 99.5
$ws.Range("L67").Value = 846.75
$ws.Range("K67").Value = 99.5
$ws.Range("M67").Value = 680.5
$ws.Range("N67").Value = -2406.75
$ws.Range("H67").Value = 597.6667
$ws.Range("J67").Value = 846.75
$ws.Range("J86").Value = 2822.25
$ws.Range("H86").Value = 40214.45
$ws.Range("I86").Value = 66608.94
$ws.Range("M86").Value = -65485.94
$ws.Range("N86").Value = -5068.25
$ws.Range("L86").Value = 2822.25
$ws.Range("K86").Value = 66608.94
$ws.Range("I89").Value = 66608.94
$ws.Range("K89").Value = 333044.7
$ws.Range("M89").Value = -327428.7
$ws.Range("N89").Value = -25343.25
$ws.Range("L89").Value = 14111.25
$ws.Range("H89").Value = 40214.45
$ws.Range("J89").Value = 2822.25
$ws.Range("M94").Value = -210
$ws.Range("N94").Value = -1834.5833
$ws.Range("L94").Value = 932.5833
$ws.Range("K94").Value = 661
$ws.Range("H94").Value = 852.7059
$ws.Range("J94").Value = 932.5833
$ws.Range("I94").Value = 661
$ws.Range("M96").Value = -5504
$ws.Range("I96").Value = 8250
$ws.Range("H96").Value = 12600
$ws.Range("K96").Value = 8250
$ws.Range("I105").Value = 204334
$ws.Range("L105").Value = 202838.2
$ws.Range("K105").Value = 204334
$ws.Range("M105").Value = -202587
$ws.Range("N105").Value = -206332.2
$ws.Range("H105").Value = 203586.1
$ws.Range("J105").Value = 202838.2
$ws.Range("I107").Value = 100040410
$ws.Range("L107").Value = 643
$ws.Range("K107").Value = 100040410
$ws.Range("M107").Value = -100038490
$ws.Range("N107").Value = -4483
$ws.Range("H107").Value = 71457620
$ws.Range("J107").Value = 643
$ws.Range("J112").Value = 36000
$ws.Range("N112").Value = -38954
$ws.Range("L112").Value = 36000
$ws.Range("H112").Value = 36000
$ws.Range("H134").Value = 2925.2812
$ws.Range("I134").Value = 2929.3225
$ws.Range("K134").Value = 8787.9675
$ws.Range("M134").Value = -6252.967500000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 33311
$ws.Range("J104").Value = 33311
$ws.Range("N104").Value = -38553
$ws.Range("L104").Value = 33311
$ws.Range("I107").Value = 9223.917
$ws.Range("L107").Value = 804
$ws.Range("K107").Value = 9223.917
$ws.Range("M107").Value = -7303.916999999999
$ws.Range("N107").Value = -4644
$ws.Range("H107").Value = 8021.0713
$ws.Range("J107").Value = 804
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K18").Value = 1004.64705
$ws.Range("M18").Value = -835.6470499999999
$ws.Range("H18").Value = 362.3
$ws.Range("I18").Value = 334.88235
$ws.Range("L39").Value = 8600.000100000001
$ws.Range("H39").Value = 2592.9285
$ws.Range("J39").Value = 2866.6667
$ws.Range("N39").Value = -9188.000100000001
$ws.Range("N43").Value = -6429.999899999999
$ws.Range("L43").Value = 6201.999899999999
$ws.Range("H43").Value = 2067.3333
$ws.Range("J43").Value = 2067.3333
$ws.Range("H99").Value = 1767.1111
$ws.Range("K99").Value = 3434.4
$ws.Range("M99").Value = -1188.4
$ws.Range("J99").Value = 2545
$ws.Range("I99").Value = 1144.8
$ws.Range("N99").Value = -12127
$ws.Range("L99").Value = 7635
$ws.Range("I107").Value = 482.77777
$ws.Range("L107").Value = 6226.200000000001
$ws.Range("K107").Value = 1448.33331
$ws.Range("M107").Value = 471.66669
$ws.Range("N107").Value = -10066.2
$ws.Range("H107").Value = 1051.5714
$ws.Range("J107").Value = 2075.4
$ws.Range("N131").Value = -12658.4505
$ws.Range("L131").Value = 2578.4505
$ws.Range("H131").Value = 824.53
$ws.Range("J131").Value = 859.4835

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K97").Value = 166669820
$ws.Range("M97").Value = -166669324
$ws.Range("N97").Value = -2741
$ws.Range("H97").Value = 125002810
$ws.Range("J97").Value = 1749
$ws.Range("I97").Value = 166669820
$ws.Range("L97").Value = 1749
$ws.Range("H104").Value = 45332
$ws.Range("J104").Value = 45332
$ws.Range("N104").Value = -52320
$ws.Range("L104").Value = 45332
$ws.Range("I107").Value = 549.7778
$ws.Range("L107").Value = 1444170.1
$ws.Range("K107").Value = 549.7778
$ws.Range("M107").Value = 1370.2222
$ws.Range("N107").Value = -1448010.1
$ws.Range("H107").Value = 632133.7
$ws.Range("J107").Value = 1444170.1
$ws.Range("H132").Value = 2038.1923
$ws.Range("I132").Value = 1527.3889
$ws.Range("M132").Value = -2052.1667
$ws.Range("K132").Value = 4582.1667

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N64").ClearContents()
$ws.Range("L64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("J96").Value = 15298.8
$ws.Range("H96").Value = 15298.8
$ws.Range("N96").Value = -20790.8
$ws.Range("L96").Value = 15298.8
$ws.Range("I100").Value = 1450
$ws.Range("L100").Value = 1794.1177
$ws.Range("K100").Value = 1450
$ws.Range("M100").Value = -909
$ws.Range("N100").Value = -2876.1177
$ws.Range("H100").Value = 1757.8948
$ws.Range("J100").Value = 1794.1177
$ws.Range("I136").Value = 2127.182
$ws.Range("L136").Value = 7048.5
$ws.Range("K136").Value = 6381.545999999999
$ws.Range("M136").Value = -3831.545999999999
$ws.Range("N136").Value = -12148.5
$ws.Range("H136").Value = 2161.3845
$ws.Range("J136").Value = 2349.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N54").Value = -7969.5
$ws.Range("L54").Value = 6929.5
$ws.Range("H54").Value = 6949.5713
$ws.Range("J54").Value = 6929.5
$ws.Range("J96").Value = 1188.5
$ws.Range("H96").Value = 71429850
$ws.Range("N96").Value = -3934.5
$ws.Range("L96").Value = 1188.5

